$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.391.45"
$ws.Range("E2").Value = "  -3.51%  "

$ws.Range("D3").Value = "3.759.17"
$ws.Range("E3").Value = "  -3.54%  "

$ws.Range("E4").Value = "  +0.02%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "613.29"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.93%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "182.63"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.85%  "

$ws.Range("D7").Value = "3.756.98"
$ws.Range("E7").Value = "  -3.46%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.638"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -4.89%  "

$ws.Range("E9").Value = "  -0.02%  "

$ws.Range("E10").Value = "  -4.11%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.165"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -8.42%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "57.31"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +2.27%  "

$ws.Range("E13").Value = "  -7.64%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "10.78"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -5.83%  "

$ws.Range("D15").Value = "4.364.72"
$ws.Range("E15").Value = "  -3.50%  "

$ws.Range("D16").Value = "3.760.40"
$ws.Range("E16").Value = "  -3.14%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "19.73"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -4.75%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "13.10"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -5.74%  "

$ws.Range("B19").Value = "Polygon"
$ws.Range("C19").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "1.14"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -5.56%  "

$ws.Range("B20").Value = "TRON"
$ws.Range("C20").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.127"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -1.99%  "

$ws.Range("D21").Value = "69.288.54"
$ws.Range("E21").Value = "  -3.34%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "417.80"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -5.09%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "4.67"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -2.07%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "89.82"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -4.79%  "

$ws.Range("E25").Value = "  -5.75%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "11.10"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -4.60%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "12.84"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -7.73%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "3.83"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -4.05%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "6.07"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +1.38%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "9.66"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -7.69%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "33.51"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -4.84%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "7.39"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -15.18%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "12.75"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -6.55%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.121"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -4.07%  "

$ws.Range("E35").Value = "  -3.40%  "

$ws.Range("B36").Value = "InjectiveProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "45.03"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -5.92%  "

$ws.Range("B37").Value = "Bittensor"
$ws.Range("C37").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "621.91"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -3.20%  "

$ws.Range("D38").Value = "0.0₃0901"
$ws.Range("E38").Value = "  -10.90%  "

$ws.Range("B39").Value = "Dai"
$ws.Range("C39").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +0.01%  "

$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.404"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -5.51%  "

$ws.Range("E41").Value = "  +0.12%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.143"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -1.95%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "3.09"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -6.65%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.0446"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -5.05%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "2.67"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -4.10%  "

$ws.Range("B46").Value = "THORChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "9.30"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -9.53%  "

$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.138"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -4.03%  "

$ws.Range("B48").Value = "dogwifhat"
$ws.Range("C48").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "2.79"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -18.60%  "

$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "2.840.95"
$ws.Range("E49").Value = "  -2.78%  "

$ws.Range("E50").Value = "  -2.30%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "3.18"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -4.37%  "
